# Apply crypto price/volume updates per commit "Updated cryptos list on Fri Apr 14 17:07:33 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.235.79'
$ws.Range("E2").Value = '  -0.43%  '
$ws.Range("D3").Value = '2.071.29'
$ws.Range("E3").Value = '  +3.28%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("E5").Value = '  +0.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5173'
$ws.Range("E7").Value = '  +1.63%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4319'
$ws.Range("E8").Value = '  +4.23%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08718'
$ws.Range("E9").Value = '  -0.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.56'
$ws.Range("E10").Value = '  +5.61%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.150'
$ws.Range("E11").Value = '  +1.50%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.06'
$ws.Range("E12").Value = '  -1.88%  '
$ws.Range("D13").Value = '2.061.12'
$ws.Range("E13").Value = '  +3.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.611'
$ws.Range("E14").Value = '  +0.71%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.621'
$ws.Range("E15").Value = '  +2.21%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  -0.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.60'
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06600'
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.65'
$ws.Range("E20").Value = '  -1.46%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.202'
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").Value = '30.264.73'
$ws.Range("E23").Value = '  -0.50%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.18'
$ws.Range("E24").Value = '  +2.05%  '
$ws.Range("E25").Value = '  +1.97%  '
$ws.Range("D26").Value = '2.303.20'
$ws.Range("E26").Value = '  +3.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.08'
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '161.69'
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.488'
$ws.Range("E29").Value = '  +3.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '130.06'
$ws.Range("E30").Value = '  -0.94%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.178'
$ws.Range("E31").Value = '  +3.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1061'
$ws.Range("E32").Value = '  +1.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.046'
$ws.Range("E33").Value = '  -0.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.838'
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.501'
$ws.Range("E35").Value = '  +11.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02550'
$ws.Range("E36").Value = '  +1.32%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.544'
$ws.Range("E37").Value = '  +5.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.390'
$ws.Range("E38").Value = '  -0.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06562'
$ws.Range("E39").Value = '  -0.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.45'
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2217'
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6628'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.234'
$ws.Range("E43").Value = '  +0.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.81'
$ws.Range("E45").Value = '  +2.19%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6262'
$ws.Range("E46").Value = '  +1.71%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.176'
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.594'
$ws.Range("E48").Value = '  -1.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.228'
$ws.Range("E49").Value = '  -2.70%  '

# Row 50: "Aave" entry replaced by "Quant" data (rank shuffled)
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '119.92'
$ws.Range("E50").Value = '  -3.45%  '

# Row 51: "Quant" entry replaced by "WEMIXTOKEN" data
$ws.Range("B51").Value = 'WEMIXTOKEN'
$ws.Range("C51").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.172'
$ws.Range("E51").Value = '  +5.94%  '
